$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header row (row 8) holding each manager's fantasy-team name, typed
#    in the same order the author entered them (Justin, Rapaka, Jayanth,
#    Sushma, Sampath, Anantha) so shared-string order matches, then merged
#    and styled (20% - Accent6, centered, bold, boxed) in column order.
# ---------------------------------------------------------------------------
$teamNames = @(
    @{cell="J8"; rng="J8:K8"; text="JUSTIN CHALLENGERS"},
    @{cell="M8"; rng="M8:N8"; text="Garuda Tejas"},
    @{cell="G8"; rng="G8:H8"; text="Jais Royal Challengers"},
    @{cell="P8"; rng="P8:Q8"; text="SUSHVIS CHOSEN ONES"},
    @{cell="S8"; rng="S8:T8"; text="GHOST RIDERS 6934"},
    @{cell="D8"; rng="D8:E8"; text="Anantha Team"}
)

foreach ($h in $teamNames) {
    $ws.Range($h.cell).Value = $h.text
}

$mergeOrder = @("D8:E8", "G8:H8", "J8:K8", "M8:N8", "P8:Q8", "S8:T8")
foreach ($r in $mergeOrder) {
    $rng = $ws.Range($r)
    $rng.Style = "20% - Accent6"
    $rng.Merge()
    $rng.HorizontalAlignment = -4108
    $rng.Borders.LineStyle = 1
    $rng.Font.Bold = $true
    $rng.Font.Size = 12
}

# ---------------------------------------------------------------------------
# 2. Contest 2 results: "DC vs KXI" (row 11) scores entered for each manager.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 20
$ws.Range("H11").Value = 60
$ws.Range("K11").Value = 80
$ws.Range("N11").Value = 100
$ws.Range("Q11").Value = 40
$ws.Range("T11").Value = 0

$wb.Application.CalculateFull()
